$wb = $excel.ActiveWorkbook

# --- Update hotel_info: insert a new "State" column between Hotel_Name and City ---
$hotel = $wb.Worksheets.Item("hotel_info")
$hotel.Range("C:C").Insert()
$hotel.Range("C1").Value = "State"
$hotel.Range("C2").Value = "Louisiana"

# --- Reorder sheet tabs: review_info first, hotel_info second ---
$review = $wb.Worksheets.Item("review_info")
$hotelRef = $wb.Worksheets.Item("hotel_info")
$review.Move($hotelRef)
